$d = $word.ActiveDocument

# --- 1) "Impre" + bookmark + "ssum" -> "Impressum" (plain merge of the two runs,
#        the stray mid-word _GoBack bookmark that split them disappears with it) ---
$found = $d.Content.Find.Execute("Impressum", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "Impressum", 2)

# --- 2) locate the "Martin: " paragraph that comes right after "Milena: Impressum " ---
$martinIdx = -1
$idx = 0
foreach ($para in $d.Paragraphs) {
    $idx = $idx + 1
    if ($para.Range.Text -eq "Martin: `r") {
        $martinIdx = $idx
    }
}

# --- 3) insert a brand new paragraph right after it for "Steffen: Logo" ---
$martinPara = $d.Paragraphs.Item($martinIdx)
$martinPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($martinIdx + 1)
# Use a trailing placeholder character so the bookmark insertion point below is
# never the very last character of the paragraph (avoids boundary weirdness),
# then strip the placeholder afterwards.
$newPara.Range.Text = "Steffen: LogoX"

$bmRange = $newPara.Range.Duplicate
$bmRange.Start = $bmRange.End - 2
$bmRange.End = $bmRange.Start
$d.Bookmarks.Add("_GoBack", $bmRange)

$newPara2 = $d.Paragraphs.Item($martinIdx + 1)
$placeholder = $newPara2.Range.Duplicate
$placeholder.Start = $placeholder.End - 2
$placeholder.End = $placeholder.Start + 1
$placeholder.Delete()
